$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 2).Value = "BTC"
$ws.Cells.Item(2, 3).Value = "Bitcoin"
$ws.Cells.Item(2, 4).Value = 27180
$ws.Cells.Item(2, 5).Value = 526852941206
$ws.Cells.Item(2, 6).Value = 8667506335
$ws.Cells.Item(2, 7).Value = 1.8148

$ws.Cells.Item(3, 2).Value = "ETH"
$ws.Cells.Item(3, 3).Value = "Ethereum"
$ws.Cells.Item(3, 4).Value = 1844.81
$ws.Cells.Item(3, 5).Value = 221888062436
$ws.Cells.Item(3, 6).Value = 5071649253
$ws.Cells.Item(3, 7).Value = 1.01581

$ws.Cells.Item(4, 2).Value = "USDT"
$ws.Cells.Item(4, 3).Value = "Tether"
$ws.Cells.Item(4, 4).Value = 1
$ws.Cells.Item(4, 5).Value = 83183082920
$ws.Cells.Item(4, 6).Value = 10097627538
$ws.Cells.Item(4, 7).Value = 0.0265

$ws.Cells.Item(5, 2).Value = "BNB"
$ws.Cells.Item(5, 3).Value = "BNB"
$ws.Cells.Item(5, 4).Value = 307.77
$ws.Cells.Item(5, 5).Value = 48593813452
$ws.Cells.Item(5, 6).Value = 352236282
$ws.Cells.Item(5, 7).Value = 0.87645

$ws.Cells.Item(6, 2).Value = "USDC"
$ws.Cells.Item(6, 3).Value = "USD Coin"
$ws.Cells.Item(6, 4).Value = 1
$ws.Cells.Item(6, 5).Value = 29065878683
$ws.Cells.Item(6, 6).Value = 2553878021
$ws.Cells.Item(6, 7).Value = 0.05116

$ws.Cells.Item(7, 2).Value = "XRP"
$ws.Cells.Item(7, 3).Value = "XRP"
$ws.Cells.Item(7, 4).Value = 0.473614
$ws.Cells.Item(7, 5).Value = 24621218448
$ws.Cells.Item(7, 6).Value = 558512875
$ws.Cells.Item(7, 7).Value = 0.11697

$ws.Cells.Item(8, 2).Value = "ADA"
$ws.Cells.Item(8, 3).Value = "Cardano"
$ws.Cells.Item(8, 4).Value = 0.377697
$ws.Cells.Item(8, 5).Value = 13211077752
$ws.Cells.Item(8, 6).Value = 187788242
$ws.Cells.Item(8, 7).Value = 3.79064

$ws.Cells.Item(9, 2).Value = "STETH"
$ws.Cells.Item(9, 3).Value = "Lido Staked Ether"
$ws.Cells.Item(9, 4).Value = 1843.19
$ws.Cells.Item(9, 5).Value = 12548372047
$ws.Cells.Item(9, 6).Value = 10740432
$ws.Cells.Item(9, 7).Value = 1.00286

$ws.Cells.Item(10, 2).Value = "DOGE"
$ws.Cells.Item(10, 3).Value = "Dogecoin"
$ws.Cells.Item(10, 4).Value = 0.072585
$ws.Cells.Item(10, 5).Value = 10128637043
$ws.Cells.Item(10, 6).Value = 227747667
$ws.Cells.Item(10, 7).Value = 1.51491

$ws.Cells.Item(11, 2).Value = "MATIC"
$ws.Cells.Item(11, 3).Value = "Polygon"
$ws.Cells.Item(11, 4).Value = 0.925192
$ws.Cells.Item(11, 5).Value = 8582552872
$ws.Cells.Item(11, 6).Value = 190630066
$ws.Cells.Item(11, 7).Value = 0.71227

$ws.Cells.Item(12, 2).Value = "SOL"
$ws.Cells.Item(12, 3).Value = "Solana"
$ws.Cells.Item(12, 4).Value = 20.54
$ws.Cells.Item(12, 5).Value = 8146994638
$ws.Cells.Item(12, 6).Value = 284877172
$ws.Cells.Item(12, 7).Value = 5.46737

$ws.Cells.Item(13, 2).Value = "TRX"
$ws.Cells.Item(13, 3).Value = "TRON"
$ws.Cells.Item(13, 4).Value = 0.077272
$ws.Cells.Item(13, 5).Value = 6978638437
$ws.Cells.Item(13, 6).Value = 298682212
$ws.Cells.Item(13, 7).Value = 1.13409

$ws.Cells.Item(14, 2).Value = "DOT"
$ws.Cells.Item(14, 3).Value = "Polkadot"
$ws.Cells.Item(14, 4).Value = 5.38
$ws.Cells.Item(14, 5).Value = 6652714960
$ws.Cells.Item(14, 6).Value = 97747265
$ws.Cells.Item(14, 7).Value = 1.75418

$ws.Cells.Item(15, 2).Value = "LTC"
$ws.Cells.Item(15, 3).Value = "Litecoin"
$ws.Cells.Item(15, 4).Value = 89.39
$ws.Cells.Item(15, 5).Value = 6527311005
$ws.Cells.Item(15, 6).Value = 490963046
$ws.Cells.Item(15, 7).Value = 1.70598

$ws.Cells.Item(16, 2).Value = "BUSD"
$ws.Cells.Item(16, 3).Value = "Binance USD"
$ws.Cells.Item(16, 4).Value = 1
$ws.Cells.Item(16, 5).Value = 5277206928
$ws.Cells.Item(16, 6).Value = 622775759
$ws.Cells.Item(16, 7).Value = 0.01033

$ws.Cells.Item(17, 2).Value = "SHIB"
$ws.Cells.Item(17, 3).Value = "Shiba Inu"
$ws.Cells.Item(17, 4).Value = 0.00000876
$ws.Cells.Item(17, 5).Value = 5159141059
$ws.Cells.Item(17, 6).Value = 119342926
$ws.Cells.Item(17, 7).Value = 2.01765

$ws.Cells.Item(18, 2).Value = "AVAX"
$ws.Cells.Item(18, 3).Value = "Avalanche"
$ws.Cells.Item(18, 4).Value = 14.51
$ws.Cells.Item(18, 5).Value = 4859707120
$ws.Cells.Item(18, 6).Value = 148658812
$ws.Cells.Item(18, 7).Value = 1.46754

$ws.Cells.Item(19, 2).Value = "DAI"
$ws.Cells.Item(19, 3).Value = "Dai"
$ws.Cells.Item(19, 4).Value = 1
$ws.Cells.Item(19, 5).Value = 4624210174
$ws.Cells.Item(19, 6).Value = 71877577
$ws.Cells.Item(19, 7).Value = 0.06035

$ws.Cells.Item(20, 2).Value = "WBTC"
$ws.Cells.Item(20, 3).Value = "Wrapped Bitcoin"
$ws.Cells.Item(20, 4).Value = 27213
$ws.Cells.Item(20, 5).Value = 4255248233
$ws.Cells.Item(20, 6).Value = 70526774
$ws.Cells.Item(20, 7).Value = 1.89795

$ws.Cells.Item(21, 2).Value = "UNI"
$ws.Cells.Item(21, 3).Value = "Uniswap"
$ws.Cells.Item(21, 4).Value = 5.06
$ws.Cells.Item(21, 5).Value = 3815039075
$ws.Cells.Item(21, 6).Value = 30997250
$ws.Cells.Item(21, 7).Value = 0.96803

$ws.Cells.Item(22, 2).Value = "LINK"
$ws.Cells.Item(22, 3).Value = "Chainlink"
$ws.Cells.Item(22, 4).Value = 6.49
$ws.Cells.Item(22, 5).Value = 3355324273
$ws.Cells.Item(22, 6).Value = 111388150
$ws.Cells.Item(22, 7).Value = 1.83646

$ws.Cells.Item(23, 2).Value = "LEO"
$ws.Cells.Item(23, 3).Value = "LEO Token"
$ws.Cells.Item(23, 4).Value = 3.5
$ws.Cells.Item(23, 5).Value = 3256793600
$ws.Cells.Item(23, 6).Value = 222266
$ws.Cells.Item(23, 7).Value = -0.14511

$ws.Cells.Item(24, 2).Value = "ATOM"
$ws.Cells.Item(24, 3).Value = "Cosmos Hub"
$ws.Cells.Item(24, 4).Value = 10.61
$ws.Cells.Item(24, 5).Value = 3103750673
$ws.Cells.Item(24, 6).Value = 64389617
$ws.Cells.Item(24, 7).Value = 0.81974

$ws.Cells.Item(25, 2).Value = "TON"
$ws.Cells.Item(25, 3).Value = "Toncoin"
$ws.Cells.Item(25, 4).Value = 1.93
$ws.Cells.Item(25, 5).Value = 2838382699
$ws.Cells.Item(25, 6).Value = 9825775
$ws.Cells.Item(25, 7).Value = 1.17334

$ws.Cells.Item(26, 2).Value = "OKB"
$ws.Cells.Item(26, 3).Value = "OKB"
$ws.Cells.Item(26, 4).Value = 46.94
$ws.Cells.Item(26, 5).Value = 2818922203
$ws.Cells.Item(26, 6).Value = 5520461
$ws.Cells.Item(26, 7).Value = 0.81464

$ws.Cells.Item(27, 2).Value = "XMR"
$ws.Cells.Item(27, 3).Value = "Monero"
$ws.Cells.Item(27, 4).Value = 153.98
$ws.Cells.Item(27, 5).Value = 2793570405
$ws.Cells.Item(27, 6).Value = 60314969
$ws.Cells.Item(27, 7).Value = 1.89621

$ws.Cells.Item(28, 2).Value = "ETC"
$ws.Cells.Item(28, 3).Value = "Ethereum Classic"
$ws.Cells.Item(28, 4).Value = 18.26
$ws.Cells.Item(28, 5).Value = 2578711360
$ws.Cells.Item(28, 6).Value = 54302180
$ws.Cells.Item(28, 7).Value = 1.51748

$ws.Cells.Item(29, 2).Value = "XLM"
$ws.Cells.Item(29, 3).Value = "Stellar"
$ws.Cells.Item(29, 4).Value = 0.088205
$ws.Cells.Item(29, 5).Value = 2364172625
$ws.Cells.Item(29, 6).Value = 28628415
$ws.Cells.Item(29, 7).Value = 0.34923

$ws.Cells.Item(30, 2).Value = "BCH"
$ws.Cells.Item(30, 3).Value = "Bitcoin Cash"
$ws.Cells.Item(30, 4).Value = 114.1
$ws.Cells.Item(30, 5).Value = 2214246094
$ws.Cells.Item(30, 6).Value = 49114846
$ws.Cells.Item(30, 7).Value = 1.12935

$ws.Cells.Item(31, 2).Value = "ICP"
$ws.Cells.Item(31, 3).Value = "Internet Computer"
$ws.Cells.Item(31, 4).Value = 4.9
$ws.Cells.Item(31, 5).Value = 2139754119
$ws.Cells.Item(31, 6).Value = 25716249
$ws.Cells.Item(31, 7).Value = -0.29036

$ws.Cells.Item(32, 2).Value = "TUSD"
$ws.Cells.Item(32, 3).Value = "TrueUSD"
$ws.Cells.Item(32, 4).Value = 1
$ws.Cells.Item(32, 5).Value = 2042403768
$ws.Cells.Item(32, 6).Value = 126444430
$ws.Cells.Item(32, 7).Value = 0.00991

$ws.Cells.Item(33, 2).Value = "FIL"
$ws.Cells.Item(33, 3).Value = "Filecoin"
$ws.Cells.Item(33, 4).Value = 4.55
$ws.Cells.Item(33, 5).Value = 1941800345
$ws.Cells.Item(33, 6).Value = 84534240
$ws.Cells.Item(33, 7).Value = 2.43392

$ws.Cells.Item(34, 2).Value = "LDO"
$ws.Cells.Item(34, 3).Value = "Lido DAO"
$ws.Cells.Item(34, 4).Value = 1.99
$ws.Cells.Item(34, 5).Value = 1753385583
$ws.Cells.Item(34, 6).Value = 25707531
$ws.Cells.Item(34, 7).Value = 0.50071

$ws.Cells.Item(35, 2).Value = "APT"
$ws.Cells.Item(35, 3).Value = "Aptos"
$ws.Cells.Item(35, 4).Value = 8.53
$ws.Cells.Item(35, 5).Value = 1694769748
$ws.Cells.Item(35, 6).Value = 52444768
$ws.Cells.Item(35, 7).Value = 4.94522

$ws.Cells.Item(36, 2).Value = "HBAR"
$ws.Cells.Item(36, 3).Value = "Hedera"
$ws.Cells.Item(36, 4).Value = 0.052486
$ws.Cells.Item(36, 5).Value = 1650616962
$ws.Cells.Item(36, 6).Value = 14114343
$ws.Cells.Item(36, 7).Value = 0.40534

$ws.Cells.Item(37, 2).Value = "ARB"
$ws.Cells.Item(37, 3).Value = "Arbitrum"
$ws.Cells.Item(37, 4).Value = 1.2
$ws.Cells.Item(37, 5).Value = 1530708979
$ws.Cells.Item(37, 6).Value = 172860255
$ws.Cells.Item(37, 7).Value = 3.57423

$ws.Cells.Item(38, 2).Value = "CRO"
$ws.Cells.Item(38, 3).Value = "Cronos"
$ws.Cells.Item(38, 4).Value = 0.060415
$ws.Cells.Item(38, 5).Value = 1527615455
$ws.Cells.Item(38, 6).Value = 4052824
$ws.Cells.Item(38, 7).Value = 0.35607

$ws.Cells.Item(39, 2).Value = "NEAR"
$ws.Cells.Item(39, 3).Value = "NEAR Protocol"
$ws.Cells.Item(39, 4).Value = 1.65
$ws.Cells.Item(39, 5).Value = 1497669006
$ws.Cells.Item(39, 6).Value = 49544388
$ws.Cells.Item(39, 7).Value = 3.15467

$ws.Cells.Item(40, 2).Value = "QNT"
$ws.Cells.Item(40, 3).Value = "Quant"
$ws.Cells.Item(40, 4).Value = 102.04
$ws.Cells.Item(40, 5).Value = 1484107837
$ws.Cells.Item(40, 6).Value = 11509495
$ws.Cells.Item(40, 7).Value = 1.12492

$ws.Cells.Item(41, 2).Value = "VET"
$ws.Cells.Item(41, 3).Value = "VeChain"
$ws.Cells.Item(41, 4).Value = 0.02026755
$ws.Cells.Item(41, 5).Value = 1473567429
$ws.Cells.Item(41, 6).Value = 56168694
$ws.Cells.Item(41, 7).Value = 4.76231

$ws.Cells.Item(42, 2).Value = "GGTKN"
$ws.Cells.Item(42, 3).Value = "GGTKN"
$ws.Cells.Item(42, 4).Value = 0.107808
$ws.Cells.Item(42, 5).Value = 1232509808
$ws.Cells.Item(42, 6).Value = 147538
$ws.Cells.Item(42, 7).Value = -5.96282

$ws.Cells.Item(43, 2).Value = "APE"
$ws.Cells.Item(43, 3).Value = "ApeCoin"
$ws.Cells.Item(43, 4).Value = 3.25
$ws.Cells.Item(43, 5).Value = 1196806518
$ws.Cells.Item(43, 6).Value = 46191481
$ws.Cells.Item(43, 7).Value = 1.23095

$ws.Cells.Item(44, 2).Value = "ALGO"
$ws.Cells.Item(44, 3).Value = "Algorand"
$ws.Cells.Item(44, 4).Value = 0.151985
$ws.Cells.Item(44, 5).Value = 1101296014
$ws.Cells.Item(44, 6).Value = 41182698
$ws.Cells.Item(44, 7).Value = 0.91753

$ws.Cells.Item(45, 2).Value = "GRT"
$ws.Cells.Item(45, 3).Value = "The Graph"
$ws.Cells.Item(45, 4).Value = 0.120578
$ws.Cells.Item(45, 5).Value = 1084465174
$ws.Cells.Item(45, 6).Value = 31214651
$ws.Cells.Item(45, 7).Value = 0.72404

$ws.Cells.Item(46, 2).Value = "USDP"
$ws.Cells.Item(46, 3).Value = "Pax Dollar"
$ws.Cells.Item(46, 4).Value = 0.999923
$ws.Cells.Item(46, 5).Value = 1019308375
$ws.Cells.Item(46, 6).Value = 25629773
$ws.Cells.Item(46, 7).Value = 0.03651

$ws.Cells.Item(47, 2).Value = "SAND"
$ws.Cells.Item(47, 3).Value = "The Sandbox"
$ws.Cells.Item(47, 4).Value = 0.548387
$ws.Cells.Item(47, 5).Value = 1011806792
$ws.Cells.Item(47, 6).Value = 153030235
$ws.Cells.Item(47, 7).Value = 5.65614

$ws.Cells.Item(48, 2).Value = "EOS"
$ws.Cells.Item(48, 3).Value = "EOS"
$ws.Cells.Item(48, 4).Value = 0.906998
$ws.Cells.Item(48, 5).Value = 1003975901
$ws.Cells.Item(48, 6).Value = 112917761
$ws.Cells.Item(48, 7).Value = 3.12462

$ws.Cells.Item(49, 2).Value = "FRAX"
$ws.Cells.Item(49, 3).Value = "Frax"
$ws.Cells.Item(49, 4).Value = 0.999883
$ws.Cells.Item(49, 5).Value = 1003904842
$ws.Cells.Item(49, 6).Value = 7688712
$ws.Cells.Item(49, 7).Value = 0.12686

$ws.Cells.Item(50, 2).Value = "EDGT"
$ws.Cells.Item(50, 3).Value = "Edgecoin"
$ws.Cells.Item(50, 4).Value = 1
$ws.Cells.Item(50, 5).Value = 1000765774
$ws.Cells.Item(50, 6).Value = 5635090
$ws.Cells.Item(50, 7).Value = 0.04031

$ws.Cells.Item(51, 2).Value = "RNDR"
$ws.Cells.Item(51, 3).Value = "Render"
$ws.Cells.Item(51, 4).Value = 2.69
$ws.Cells.Item(51, 5).Value = 985768314
$ws.Cells.Item(51, 6).Value = 87786009
$ws.Cells.Item(51, 7).Value = -1.23651

